# Apply "Penality Reward System" edits to B0CVHD8PLP_po_data workbook
$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Weekly Quantity" ---
# Remove the two rows that contained the weeks ending 45361.99999999999 (qty 280)
# and 45368.99999999999 (qty 100). Deleting sheet rows 4 and 5 shifts the rows
# below upward, which matches the new A1:B19 layout described in the diff.
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws1.Rows.Item(4).Delete()
$ws1.Rows.Item(4).Delete()

# --- Sheet 2: "Monthly Trend" ---
# Row 3 (week of 45382.99999999999) requested quantity changes from 480 to 100.
$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Range("B3").Value = 100
